{"js": "// The two \"Deliverable/Feature\" and \"Story\" bullets under the gold-price\n// Impact/Capability branch are being tightened up:\n//   \"Display gold prices from various jewelers\" -> \"Display gold price from various jewelers\"\n//   \"Latest prices from each website\"            -> \"View latest price\"\n\nconst body = context.document.body;\n\n// 1) \"Display gold prices from various jewelers\" -> \"Display gold price from various jewelers\"\nconst displaySearch = body.search(\"Display gold prices from various jewelers\", { matchCase: true });\ndisplaySearch.load(\"text\");\nawait context.sync();\n\nif (displaySearch.items.length > 0) {\n  displaySearch.items[0].insertText(\n    \"Display gold price from various jewelers\",\n    Word.InsertLocation.replace\n  );\n}\n\n// 2) \"Latest prices from each website\" -> \"View latest price\"\nconst storySearch = body.search(\"Latest prices from each website\", { matchCase: true });\nstorySearch.load(\"text\");\nawait context.sync();\n\nif (storySearch.items.length > 0) {\n  storySearch.items[0].insertText(\"View latest price\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The two \"Deliverable/Feature\" and \"Story\" bullets under the gold-price\n# Impact/Capability branch are being tightened up:\n#   \"Display gold prices from various jewelers\" -> \"Display gold price from various jewelers\"\n#   \"Latest prices from each website\"            -> \"View latest price\"\n\n$d = $word.ActiveDocument\n\n# wdReplace constants used below: wdFindContinue = 1, wdReplaceAll = 2\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n# 1) \"Display gold prices from various jewelers\" -> \"Display gold price from various jewelers\"\n$find1 = $d.Content.Find\n$find1.Text = \"Display gold prices from various jewelers\"\n$find1.Replacement.Text = \"Display gold price from various jewelers\"\n$find1.Execute(\n    $find1.Text, $false, $false, $false, $false, $false, $true,\n    $wdFindContinue, $false, $find1.Replacement.Text, $wdReplaceAll\n) | Out-Null\n\n# 2) \"Latest prices from each website\" -> \"View latest price\"\n$find2 = $d.Content.Find\n$find2.Text = \"Latest prices from each website\"\n$find2.Replacement.Text = \"View latest price\"\n$find2.Execute(\n    $find2.Text, $false, $false, $false, $false, $false, $true,\n    $wdFindContinue, $false, $find2.Replacement.Text, $wdReplaceAll\n) | Out-Null\n"}
